$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.945.80"
$ws.Range("E2").Value = "  +8.26%  "
$ws.Range("D3").Value = "1.814.72"
$ws.Range("E3").Value = "  +5.13%  "
$ws.Range("D4").Value = "'0.9995"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "'246.73"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.61%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("D7").Value = "'0.4932"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +2.00%  "
$ws.Range("B8").Value = "OKB"
$ws.Range("C8").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D8").Value = "'43.66"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +5.68%  "
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "'0.2780"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +7.55%  "
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "'0.06404"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +3.57%  "
$ws.Range("B11").Value = "WrappedEther"
$ws.Range("C11").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D11").Value = "1.809.23"
$ws.Range("E11").Value = "  +4.81%  "
$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").Value = "'16.75"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +5.53%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "'0.07073"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.90%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "'0.6445"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +6.68%  "
$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").Value = "'83.96"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +9.08%  "
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").Value = "'4.679"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +4.74%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "28.964.11"
$ws.Range("E17").Value = "  +9.06%  "
$ws.Range("B18").Value = "Dai"
$ws.Range("C18").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D18").Value = "'0.9985"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "'0.000007325"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.37%  "
$ws.Range("B20").Value = "BinanceUSD"
$ws.Range("C20").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D20").Value = "'1.001"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.28%  "
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Value = "'12.27"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +7.99%  "
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.042.89"
$ws.Range("E22").Value = "  +4.74%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "'4.579"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +3.71%  "
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").Value = "'8.808"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +3.10%  "
$ws.Range("B25").Value = "Chainlink"
$ws.Range("C25").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D25").Value = "'5.345"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +5.70%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "'143.08"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.08%  "
$ws.Range("B27").Value = "BitcoinCash"
$ws.Range("C27").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D27").Value = "'129.19"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +21.47%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'16.43"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +7.86%  "
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").Value = "'1.887"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +6.14%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'1.410"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +3.07%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "'4.136"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +3.04%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "'0.08349"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +5.28%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'3.779"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.92%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.04944"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +9.50%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "'1.096"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +9.58%  "
$ws.Range("D36").Value = "'2.696"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.80%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'0.6724"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +8.69%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "'2.284"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +14.07%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "'2.745"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +11.93%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "'0.9550"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.90%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'6.139"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +9.41%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "'0.01586"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +5.85%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").Value = "'1.000"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.25%  "
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").Value = "'100.87"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.09%  "
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").Value = "'0.4083"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +6.60%  "
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").Value = "'7.139"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +5.06%  "
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").Value = "'0.1221"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +5.71%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "'0.05529"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +3.17%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'8.136"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.98%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "'31.62"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +4.99%  "
$ws.Range("D51").Value = "'0.3619"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +7.93%  "
